# Update analysis with december report data.
# Column C (count "n") values for rows 2-54 are revised upward: each value
# is scaled from the November report (divisor 7) to the December report
# (divisor 8), i.e. new = old * 8 / 7 (old values are always multiples of 7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column C (data starts at row 2).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 2 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $oldVal = $cell.Value2
    if ($oldVal -ne $null) {
        $cell.Value2 = [Math]::Round(($oldVal * 8 / 7), 0)
    }
}
